$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 14
$ws.Cells.Item($row, 1).Value = 42620.886458333334
$ws.Cells.Item($row, 2).Value = -28
$ws.Cells.Item($row, 3).Value = 51
$ws.Cells.Item($row, 4).Value = 47
$ws.Cells.Item($row, 5).Value = 51
$ws.Cells.Item($row, 6).Value = 78
$ws.Cells.Item($row, 7).Value = 16853
$ws.Cells.Item($row, 8).Value = 11458
$ws.Cells.Item($row, 9).Value = 598
$ws.Cells.Item($row, 10).Value = 90
$ws.Cells.Item($row, 11).Value = 83
$ws.Cells.Item($row, 12).Value = 3
$ws.Cells.Item($row, 13).Value = 11
$ws.Cells.Item($row, 14).Value = "Named"
